# corrected data cleaning for pre/post/total fixation data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Strip the bold/border/center-align header style from row 1 ---
# (the header row no longer gets the special "bold + thin border + centered"
# formatting; every cell reverts to the default "Normal" style)
$ws.Range("A1:S1").Style = "Normal"

# A1 used to hold the literal text "Unnamed: 0" - that label goes away.
$ws.Range("A1").Value = ""

# --- 2. Corrected numeric values -------------------------------------

# Row 3 - Revisit count
$ws.Range("D3").Value = 30
$ws.Range("E3").Value = 20
$ws.Range("O3").Value = 24

# Row 4 - Fixation count
$ws.Range("D4").Value = 79
$ws.Range("E4").Value = 40
$ws.Range("O4").Value = 196

# Row 5 - Dwell time (ms)
$ws.Range("D5").Value = 26017.15
$ws.Range("E5").Value = 14205.09
$ws.Range("O5").Value = 71757.98

# Row 6 - Dwell time (%)
$ws.Range("B6").Value = 0.77
$ws.Range("C6").Value = 5.83
$ws.Range("D6").Value = 21.74
$ws.Range("E6").Value = 11.87
$ws.Range("F6").Value = 1.17
$ws.Range("H6").Value = 1.17
$ws.Range("J6").Value = 4.31
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 2.32
$ws.Range("M6").Value = 1.17
$ws.Range("O6").Value = 59.97
$ws.Range("Q6").Value = 0.28

# Row 7 - Fixation duration (ms)
$ws.Range("D7").Value = 329.33
$ws.Range("E7").Value = 355.13
$ws.Range("O7").Value = 366.11

# --- 3. Row 10 gained a row label -------------------------------------
$ws.Range("A10").Value = "pos"

# --- 4. Append a brand-new "pos" table in rows 11-18 ------------------

# Row 11: column headers (same schema as row 1, but unstyled)
$ws.Range("B11").Value = "arg"
$ws.Range("C11").Value = "arg2"
$ws.Range("D11").Value = "code"
$ws.Range("E11").Value = "condbody"
$ws.Range("F11").Value = "conditionalstate"
$ws.Range("G11").Value = "gemini"
$ws.Range("H11").Value = "literal"
$ws.Range("I11").Value = "literal2"
$ws.Range("J11").Value = "literal3"
$ws.Range("K11").Value = "methodcall"
$ws.Range("L11").Value = "methodcall2"
$ws.Range("M11").Value = "methoddec"
$ws.Range("N11").Value = "param"
$ws.Range("O11").Value = "summary"
$ws.Range("P11").Value = "var"
$ws.Range("Q11").Value = "var2"
$ws.Range("R11").Value = "var3"
$ws.Range("S11").Value = "var5"

# Row 12: section header
$ws.Range("A12").Value = "Fixation based metrics"

# Row 13: Revisit count
$ws.Range("A13").Value = "Revisit count"
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 0
$ws.Range("G13").Value = 7
$ws.Range("O13").Value = 10

# Row 14: Fixation count
$ws.Range("A14").Value = "Fixation count"
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 1
$ws.Range("G14").Value = 123
$ws.Range("O14").Value = 76

# Row 15: Dwell time (ms)
$ws.Range("A15").Value = "Dwell time (ms)"
$ws.Range("D15").Value = 2360.8
$ws.Range("E15").Value = 183.5
$ws.Range("G15").Value = 35234.67
$ws.Range("O15").Value = 20153.75

# Row 16: Dwell time (%)
$ws.Range("A16").Value = "Dwell time (%)"
$ws.Range("D16").Value = 2.58
$ws.Range("E16").Value = 0.2
$ws.Range("G16").Value = 38.55
$ws.Range("O16").Value = 22.05

# Row 17: Fixation duration (ms)
$ws.Range("A17").Value = "Fixation duration (ms)"
$ws.Range("D17").Value = 472.16
$ws.Range("E17").Value = 183.5
$ws.Range("G17").Value = 286.46
$ws.Range("O17").Value = 265.18

# Row 18: First fixation duration (ms)
$ws.Range("A18").Value = "First fixation duration (ms)"
$ws.Range("D18").Value = 1234.56
$ws.Range("E18").Value = 183.5
$ws.Range("G18").Value = 250.16
$ws.Range("O18").Value = 216.86
